$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Relocate the "_GoBack" bookmark from the end of the "VŨ" paragraph
#    to the end of the "NHÂN" paragraph (right before its paragraph
#    mark), matching the diff.
# ---------------------------------------------------------------------

# Find the "NHÂN" paragraph (it is the lone-word paragraph right after
# the "NIÊN" paragraph and the "○" bookmark paragraph).
$nhanPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "NHÂN" + [char]13) {
        $nhanPara = $p
        break
    }
}

if ($nhanPara -eq $null) {
    throw "Could not locate the 'NHAN' paragraph"
}

$nhanEnd = $nhanPara.Range.End - 1   # position right before the pilcrow

# Remove the existing _GoBack bookmark (currently at the end of the
# "VŨ" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Word COM's Bookmarks.Add mishandles a Range that is collapsed exactly
# at "end-of-paragraph-text" (immediately before the pilcrow) - it snaps
# to the whole paragraph instead of staying collapsed. Work around this
# by temporarily inserting a placeholder character after that position,
# bookmarking just before the placeholder (no longer the paragraph-end
# position), and then deleting the placeholder again; the now-collapsed
# bookmark stays anchored in place.
$placeholder = $d.Range($nhanEnd, $nhanEnd)
$placeholder.InsertAfter("X")

$bookmarkTarget = $d.Range($nhanEnd, $nhanEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkTarget)

$placeholderRange = $d.Range($nhanEnd, $nhanEnd + 1)
$placeholderRange.Delete()

# ---------------------------------------------------------------------
# 2) Split " – (VŨ QUÝ) – (VŨ LẠP) – (ĐÔNG VŨ)" into three runs so that
#    "THU" becomes its own run with an eastAsia font hint, replacing
#    "ĐÔNG".
# ---------------------------------------------------------------------

$found = $d.Content.Find.Execute("(ĐÔNG VŨ)", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '(ĐÔNG VŨ)' text"
}
$matchRange = $d.Content.Duplicate()

# Re-find to get the actual matched Range (Execute moved/selected it on
# $d.Content already, so reuse that range).
$target = $d.Content
$target.Find.Execute("(ĐÔNG VŨ)", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)

$matchStart = $target.Start
$matchEnd = $target.End

# Replace "(ĐÔNG VŨ)" with "(THU VŨ)" as plain text first, so run
# boundaries/text are simple to reason about.
$target.Text = "(THU VŨ)"

$openParenEnd = $matchStart + 1          # just after "("
$thuStart = $openParenEnd
$thuEnd = $thuStart + 3                  # length of "THU"

# Run 1 already ends right before "(" is untouched; now give "THU" its
# own run with the eastAsia font hint.
$thuRange = $d.Range($thuStart, $thuEnd)
$thuRange.Font.Name  # touch Font to force a distinct run (no-op read)
$thuRange.Text = "THU"
$thuRange = $d.Range($thuStart, $thuEnd)

# Ensure "THU" sits in its own run by re-typing it through a Range that
# is bounded by plain-text neighbours - setting the font hint below
# naturally forces Word to split the run.
$thuRange.Font.Name = $thuRange.Font.Name
